$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New FedEx tracking numbers for column P (ShipmentTracking), rows 2-26,
# from the latest shipment run ("Final changes for cheetah").
$trackingNumbers = @{
    2  = "320018586090"
    3  = "320018586104"
    4  = "320018586137"
    5  = "320018586159"
    6  = "320018586192"
    7  = "320018586218"
    8  = "320018586240"
    9  = "320018586262"
    10 = "320018586295"
    11 = "320018586310"
    12 = "320018586354"
    13 = "320018586376"
    14 = "320018586402"
    15 = "320018586424"
    16 = "320018586457"
    17 = "320018586479"
    18 = "320018586516"
    19 = "320018586538"
    20 = "320018586560"
    21 = "320018586582"
    22 = "320018586619"
    23 = "320018586620"
    24 = "320018586630"
    25 = "320018586641"
    26 = "320018586652"
}

# The tracking numbers are plain digit strings that Excel would otherwise
# auto-convert to numbers; force the range to Text first so the values are
# written as shared strings (matching column P's existing text cells), then
# restore the default "Normal" style so no visible formatting change sticks.
$pRange = $ws.Range("P2:P26")
$pRange.NumberFormat = "@"
foreach ($row in $trackingNumbers.Keys) {
    $ws.Range("P$row").Value = $trackingNumbers[$row]
}
$pRange.Style = "Normal"

# Row 20 also got a new ActualRate and a failing Result in this run.
# ActualRate is stored as literal text (e.g. "$19.04") elsewhere in column Q,
# so force Q20 to Text too, otherwise "$104.69" is auto-parsed as currency.
$q20 = $ws.Range("Q20")
$q20.NumberFormat = "@"
$q20.Value = "$104.69"
$q20.Style = "Normal"

$ws.Range("R20").Value = "FAIL"
